$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) The contract text itself barely changes: the phrase "в однокомнатных и
#    двухкомнатных номерах" loses " и двухкомнатных", becoming simply
#    "в однокомнатных номерах".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "однокомнатных и двухкомнатных номерах", $true, $false, $false, $false,
    $false, $true, 1, $false, "однокомнатных номерах", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) The `_GoBack` bookmark (Word's "last edit position" marker) moves from
#    right after ", Москва" earlier in the document to right after
#    "однокомнатных" in the sentence we just edited. Word only ever keeps a
#    single `_GoBack` bookmark, so re-adding it under that name automatically
#    relocates (and the remaining bookmark ids renumber accordingly).
# ---------------------------------------------------------------------------
$target = $d.Content
$target.Find.Execute(
    "однокомнатных", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null
$target.Collapse(0)   # wdCollapseEnd
$d.Bookmarks.Add("_GoBack", $target) | Out-Null

# ---------------------------------------------------------------------------
# 3) The footer's cached PAGE field result text changes from "1" to "2".
# ---------------------------------------------------------------------------
$footer = $d.Sections.Item(1).Footers.Item(1)
$footer.Range.Find.Execute(
    "1", $true, $false, $false, $false, $false,
    $true, 1, $false, "2", 2) | Out-Null
